$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.204.83"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "3.165.56"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("D5").Value = "'577.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.83%  "
$ws.Range("D6").Value = "'150.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.31%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.164.64"
$ws.Range("E8").Value = "  +3.67%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.24%  "
$ws.Range("D11").Value = "'6.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "'0.501"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.50%  "
$ws.Range("E13").Value = "  +16.62%  "
$ws.Range("D14").Value = "'37.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.39%  "
$ws.Range("D15").Value = "3.683.23"
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("D16").Value = "65.253.97"
$ws.Range("E16").Value = "  +2.37%  "
$ws.Range("D17").Value = "3.162.07"
$ws.Range("E17").Value = "  +3.92%  "
$ws.Range("D18").Value = "'7.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.80%  "
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").Value = "'511.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.07%  "
$ws.Range("E21").Value = "  +4.54%  "
$ws.Range("D22").Value = "'0.723"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.60%  "
$ws.Range("D23").Value = "'15.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.64%  "
$ws.Range("D24").Value = "'7.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.96%  "
$ws.Range("D25").Value = "'84.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'9.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.11%  "
$ws.Range("D28").Value = "'2.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.10%  "
$ws.Range("E29").Value = "  +7.91%  "
$ws.Range("E30").Value = "  +15.00%  "
$ws.Range("D31").Value = "'27.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.28%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("D34").Value = "'6.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.48%  "
$ws.Range("D35").Value = "'6.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.26%  "
$ws.Range("D36").Value = "'55.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "'0.0904"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.25%  "
$ws.Range("D38").Value = "'472.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.06%  "
$ws.Range("D39").Value = "'0.0421"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("D40").Value = "'3.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.39%  "
$ws.Range("D41").Value = "'8.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.19%  "
$ws.Range("D42").Value = "3.061.86"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("E44").Value = "  +11.58%  "
$ws.Range("E45").Value = "  +5.15%  "
$ws.Range("D46").Value = "'28.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.88%  "
$ws.Range("D47").Value = "0.0₃0599"
$ws.Range("E47").Value = "  +16.73%  "
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("E50").Value = "  +7.70%  "
$ws.Range("D51").Value = "'121.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.73%  "

Write-Host "Updated cryptos list"